# The workbook contains 16 worksheets (one per backward-elimination step),
# each with the full statsmodels OLS summary text dumped into cell B2.
# The summary text embeds the timestamp the regression was run at; the
# commit re-ran the regressions a few days later, so every sheet's "Date:"
# and "Time:" fields need to be updated to the new run's timestamp while
# leaving the rest of the (already-recomputed) summary text untouched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2
    if ($text -ne $null -and $text -like "*Date:*") {
        $text = $text.Replace("Thu, 02 Jan 2020", "Sun, 05 Jan 2020")
        $text = $text.Replace("20:48:45", "21:22:23")
        $cell.Value2 = $text
    }
}
